# Insert a new data row at row 32 (pushing the existing rows 32..68 down to
# 33..69) and populate it with a new price report entry, as described by the
# commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 32 (shifts rows 32-68 down to 33-69)
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record
$ws.Cells.Item(32, 1).Value  = 1
$ws.Cells.Item(32, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value  = 44902
$ws.Cells.Item(32, 5).Value  = 15
$ws.Cells.Item(32, 6).Value  = "Fruta"
$ws.Cells.Item(32, 7).Value  = 100103
$ws.Cells.Item(32, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(32, 9).Value  = 100103006
$ws.Cells.Item(32, 10).Value = "Nectarín"
$ws.Cells.Item(32, 11).Value = "Super Queen"
$ws.Cells.Item(32, 12).Value = "Segunda"
$ws.Cells.Item(32, 13).Value = 400
$ws.Cells.Item(32, 14).Value = 18000
$ws.Cells.Item(32, 15).Value = 19000
$ws.Cells.Item(32, 16).Value = 18500
$ws.Cells.Item(32, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(32, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(32, 19).Value = 1028
$ws.Cells.Item(32, 20).Value = 18
